# issue #5: stock data from json to db
#
# The "股票" (stock) sheet (4th worksheet) gains three new columns that
# mirror the JSON-to-DB export shape used for the other property sheets:
#   - "category"    -> literal "normal"   (inserted right after property_category)
#   - "source_file" -> literal "tmpbc031" (appended after legislator_id)
#   - "index"        -> same value as the existing row index in column A

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# --- Insert the new "category" column between "property_category" (H) and
#     "date" (old I). This shifts date/legislator_name/legislator_id one
#     column to the right (I->J, J->K, K->L) and carries their cell styles
#     along automatically.
$ws.Columns("I").Insert()

$ws.Range("I1").Value2 = "category"

# --- Append the two trailing columns: source_file, index.
$ws.Range("M1").Value2 = "source_file"
$ws.Range("N1").Value2 = "index"

# Match the header formatting (bold + border, style used by the rest of row 1).
$ws.Range("L1").Copy()
$ws.Range("M1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$lastRow = 11
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("I" + $r).Value2 = "normal"
    $ws.Range("M" + $r).Value2 = "tmpbc031"
    $ws.Range("N" + $r).Value2 = $ws.Range("A" + $r).Value2
}
